$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '22.406.70'
$ws.Range("E2").Value = '  +0.01%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.570.74'
$ws.Range("E3").Value = '  -0.10%  '
$ws.Range("E4").Value = '  +0.15%  '
$ws.Range("E5").Value = '  +0.11%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '291.49'
$ws.Range("E6").Value = '  +0.24%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.3732'
$ws.Range("E7").Value = '  -1.00%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '49.81'
$ws.Range("E8").Value = '  -0.11%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.3389'
$ws.Range("E9").Value = '  -0.86%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.07540'
$ws.Range("E10").Value = '  -1.45%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '1.132'
$ws.Range("E11").Value = '  -2.31%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.003'
$ws.Range("E12").Value = '  -0.02%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '21.33'
$ws.Range("E13").Value = '  +0.31%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.983'
$ws.Range("E14").Value = '  -0.52%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '6.919'
$ws.Range("E15").Value = '  -0.24%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '1.578.39'
$ws.Range("E16").Value = '  +0.38%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.00001118'
$ws.Range("E17").Value = '  -1.38%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '90.88'
$ws.Range("E18").Value = '  +1.03%  '
$ws.Range("E19").Value = '  -0.29%  '
$ws.Range("E20").Value = '  +0.10%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.259'
$ws.Range("E21").Value = '  +0.76%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '16.32'
$ws.Range("E22").Value = '  -2.99%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '12.10'
$ws.Range("E23").Value = '  +0.58%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '22.415.08'
$ws.Range("E24").Value = '  +0.10%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.331'
$ws.Range("E25").Value = '  -3.65%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.623'
$ws.Range("E26").Value = '  -3.54%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '20.10'
$ws.Range("E27").Value = '  -0.76%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '148.41'
$ws.Range("E28").Value = '  +1.20%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '5.014'
$ws.Range("E29").Value = '  -0.33%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '125.46'
$ws.Range("E30").Value = '  -0.64%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.752.97'
$ws.Range("E31").Value = '  +0.56%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.045'
$ws.Range("E33").Value = '  -1.13%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.965'
$ws.Range("E34").Value = '  -2.22%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '9.728'
$ws.Range("E35").Value = '  -2.99%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.08363'
$ws.Range("E36").Value = '  -2.78%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.378'
$ws.Range("E37").Value = '  +3.57%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.02459'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.2284'
$ws.Range("E39").Value = '  -1.35%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.06495'
$ws.Range("E40").Value = '  -1.29%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '5.435'
$ws.Range("E41").Value = '  -0.61%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '11.24'
$ws.Range("E42").Value = '  -2.61%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.6193'
$ws.Range("E43").Value = '  -3.96%  '
$ws.Range("E44").Value = '  +0.01%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '13.83'
$ws.Range("E45").Value = '  -1.91%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '3.812'
$ws.Range("E46").Value = '  +0.42%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.5783'
$ws.Range("E47").Value = '  -3.88%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '129.29'
$ws.Range("E48").Value = '  +3.15%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.068'
$ws.Range("E49").Value = '  -0.81%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.211'
$ws.Range("E50").Value = '  -7.01%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.07310'
$ws.Range("E51").Value = '  -0.23%  '
